# Apply crypto price/volume updates per commit "Updated cryptos list on Mon May 13 11:47:45 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.641.34"
$ws.Range("E2").Value = "  +2.44%  "
$ws.Range("D3").Value = "2.961.32"
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "2.959.98"
$ws.Range("E8").Value = "  +1.18%  "
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.85%  "
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.448"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").Value = "3.449.25"
$ws.Range("D17").Value = "62.524.00"
$ws.Range("E17").Value = "  +2.43%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "2.951.34"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "441.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("E26").Value = "  +0.86%  "
$ws.Range("E27").Value = "  -3.80%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.64"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("E33").Value = "  -2.25%  "
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").Value = "0.0₃0876"
$ws.Range("E35").Value = "  +1.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.993"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.64"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.62"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.57%  "
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.281"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.26%  "
$ws.Range("D45").Value = "2.721.66"
$ws.Range("E45").Value = "  +1.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "136.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.86%  "
$ws.Range("E47").Value = "  -1.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "363.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.37%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.87%  "
